# Auto-generated edit script: updates cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.797.80"
$ws.Range("D3").Value = "'1.868.90"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'300.28"
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.5367"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").Value = "'0.07144"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "'21.53"
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("D11").Value = "'0.8872"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "'0.08137"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "'1.912.95"
$ws.Range("E13").Value = "  +47.87%  "
$ws.Range("D14").Value = "'92.41"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").Value = "'5.287"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'14.85"
$ws.Range("D18").Value = "'0.000008482"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'26.837.62"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'4.966"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "'10.65"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "'6.376"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "'2.286"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'146.21"
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("D26").Value = "'1.741"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D28").Value = "'113.70"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "'4.697"
$ws.Range("E29").Value = "  -2.77%  "
$ws.Range("D30").Value = "'4.625"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").Value = "'0.8113"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").Value = "'0.05018"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").Value = "'1.172"
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("D35").Value = "'2.948"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").Value = "'0.6062"
$ws.Range("E36").Value = "  +5.31%  "
$ws.Range("D37").Value = "'2.678"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'3.196"
$ws.Range("E38").Value = "  -4.83%  "
$ws.Range("D39").Value = "'0.01946"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").Value = "'0.5280"
$ws.Range("E41").Value = "  +7.28%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.476"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'8.761"
$ws.Range("E43").Value = "  -6.81%  "
$ws.Range("D44").Value = "'116.29"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "'0.1488"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").Value = "'1.643"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "'37.27"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("D50").Value = "'0.06055"
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("E51").Value = "  -2.59%  "
